$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
try {
  $sm.Theme.Name = "TestName"
  Write-Output "set name ok"
} catch {
  Write-Output "err: $_"
}
Write-Output "Name now = [$($sm.Theme.Name)]"
